# Update final-project due date and add an "Using LLMs Well" reading link.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) The "Plotting" week reading list (row 24, column C) gains a new bullet
#    pointing to the "Using LLMs Well" notebook.
$plottingCell = $ws.Range("C24")
$currentReadings = $plottingCell.Value2
$newBullet = "- ``Using LLMs Well <../notebooks/PDS_not_yet_in_coursera/99_advice/using_llms_well.html>```_"
$plottingCell.Value = $currentReadings + "`n" + $newBullet

# 2) The final project due-date cell (row 28, column A) moves from
#    "Wed Dec 12" to "Wed Dec 10".
$ws.Range("A28").Value = "Wed Dec 10"

# 3) Reflect the author's updated on-screen selection/scroll position when the
#    file was saved (best effort - scroll position, selection landed on B24).
$ws.Activate()
$ws.Range("A22").Select()
$ws.Range("B24").Select()
